# Apply "add colors TOPAZE and AMETHYSTE" edit to schedule_retour.xlsx
#
# Summary of the change:
#  - Header row fill becomes a medium purple ("AMETHYSTE", #9966CC) and its
#    (bold, white) font gains an explicit Arial family name.
#  - The rows that used to be light blue (#BDD7EE) become lavender (#E6E6FA).
#  - The rows that used to be light green (#E2EFDA) become thistle
#    ("TOPAZE", #D8BFD8) - a brand new color added to the palette.
#  - Column F gets one character narrower (10 -> 9).
#  - The header labels are re-worded / translated to French.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Convert-HexToBgr {
    param([string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$amethyste = Convert-HexToBgr "9966CC"   # new header fill color
$lavande   = Convert-HexToBgr "E6E6FA"   # new color for former light-blue rows
$topaze    = Convert-HexToBgr "D8BFD8"   # brand new color for former light-green rows

# ---------------------------------------------------------------------
# 1) Header row (row 1): recolor fill to AMETHYSTE and add Arial to the
#    existing bold/white font, then update the header captions.
# ---------------------------------------------------------------------
$header = $ws.Range("A1:G1")
$header.Interior.Color = $amethyste
$header.Interior.PatternColor = $amethyste
$header.Font.Name = "Arial"

$ws.Range("A1").Value = "Round"
$ws.Range("B1").Value = "Début"
$ws.Range("C1").Value = "Fin"
$ws.Range("D1").Value = "Équipe 1"
$ws.Range("E1").Value = "Équipe 2"
$ws.Range("F1").Value = "Durée"
$ws.Range("G1").Value = "Phase"

# ---------------------------------------------------------------------
# 2) Rows that were light blue (#BDD7EE) -> lavender (#E6E6FA)
# ---------------------------------------------------------------------
$lavenderRows = @(2,4,6,8,10,12,14,15,17,18,19,20,21,22,23,24,25,26,27,28,29)
$lavenderAddr = ($lavenderRows | ForEach-Object { "A${_}:G${_}" }) -join ","
$lavenderRange = $ws.Range($lavenderAddr)
foreach ($area in $lavenderRange.Areas) {
    $area.Interior.Color = $lavande
    $area.Interior.PatternColor = $lavande
}

# ---------------------------------------------------------------------
# 3) Rows that were light green (#E2EFDA) -> TOPAZE (#D8BFD8, brand new)
# ---------------------------------------------------------------------
$topazeRows = @(3,5,7,9,11,13,16)
$topazeAddr = ($topazeRows | ForEach-Object { "A${_}:G${_}" }) -join ","
$topazeRange = $ws.Range($topazeAddr)
foreach ($area in $topazeRange.Areas) {
    $area.Interior.Color = $topaze
    $area.Interior.PatternColor = $topaze
}

# ---------------------------------------------------------------------
# 4) Column F is one character narrower now (10 -> 9)
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 8.14
